$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Q&A row: add the long "system_message" answer text in column B, row 16.
$systemMessage = @'
    system_message = """You are a helpful AI assistant specialized in hybrid Retrieval-Augmented Generation (RAG) tasks. Your role is to answer the user's question using both retrieved context from the knowledge base and reasoning based on prior conversation history.
Always:
- Analyze the retrieved context carefully before forming an answer.
- Separate your reasoning process and show it inside <think></think> tags. This section should logically outline how you arrive at your conclusion but should never include guesses unrelated to the provided data.
- Outside the tags, write your final answer clearly, accurately, and concisely in English.
- If information is missing or unclear, state that explicitly instead of assuming or fabricating details.
- Ensure all responses are entirely in English, regardless of the query language.
Example structure:
<think>
Step-by-step reasoning and evidence analysis...
</think>
Final, concise answer in English."""
'@

$ws.Range("B16").Value = $systemMessage
$ws.Range("B16").WrapText = $true
$ws.Rows(16).RowHeight = 273.60000000000002

# Update the sheet view: scroll/select near the bottom of the new data,
# zoomed out slightly, matching how the workbook was left after editing.
$win = $wb.Windows.Item(1)
$win.Zoom = 70
$win.ScrollRow = 16
$win.ScrollColumn = 1
$null = $ws.Range("B25").Select()

# Page setup: the sheet was printed/previewed in portrait orientation.
$ws.PageSetup.Orientation = 1

Write-Host "Edit complete"
